$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Create the new shared strings in exactly the order needed so the
#    resulting sharedStrings table lines up with the target (indices
#    8..14 appended after the 8 existing ones).
# ------------------------------------------------------------------
$ws.Cells.Item(2,3).Value  = "거래일"
$ws.Cells.Item(2,1).Value  = "논리명"
$ws.Cells.Item(3,1).Value  = "물리명"
$ws.Cells.Item(4,1).Value  = "인포타입"
$ws.Cells.Item(5,1).Value  = "데이터타입"
$ws.Cells.Item(1,1).Value  = "1. 일간시세"
$ws.Cells.Item(7,1).Value  = "2. 관심종목"

# Reuse the label strings for the second ("관심종목") block.
$ws.Cells.Item(8,1).Value  = "논리명"
$ws.Cells.Item(9,1).Value  = "물리명"
$ws.Cells.Item(10,1).Value = "인포타입"
$ws.Cells.Item(11,1).Value = "데이터타입"

# ------------------------------------------------------------------
# 2) Move the original header row down into row 2, inserting a new
#    "거래일" column between 종목코드 and 종가.
# ------------------------------------------------------------------
$ws.Cells.Item(2,2).Value  = "종목코드"
$ws.Cells.Item(2,4).Value  = "종가"
$ws.Cells.Item(2,5).Value  = "시가"
$ws.Cells.Item(2,6).Value  = "고가"
$ws.Cells.Item(2,7).Value  = "저가"
$ws.Cells.Item(2,8).Value  = "거래량"
$ws.Cells.Item(2,9).Value  = "기관순매매"
$ws.Cells.Item(2,10).Value = "외국인순매매"

# The old row 1 (B1:H1) content has now been relocated to row 2 -
# clear what remains of the old header row except A1.
$ws.Range("B1:H1").ClearContents()

# ------------------------------------------------------------------
# 3) Formatting.
# ------------------------------------------------------------------
# Label cells (A2:A5, A8:A11): solid theme "Accent 6" fill, regular font.
foreach ($addr in @("A2","A3","A4","A5","A8","A9","A10","A11")) {
    $c = $ws.Range($addr)
    $c.Interior.Pattern = 1
    $c.Interior.ThemeColor = 10
    $c.Interior.TintAndShade = 0
}

# Section heading "2. 관심종목" (A7): bold.
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Interior.ColorIndex = -4142

# Section heading "1. 일간시세" (A1): bold.
$ws.Range("A1").Font.Bold = $true

# ------------------------------------------------------------------
# 4) Selection / active cell.
# ------------------------------------------------------------------
$ws.Range("B10").Select()
